# Insert a new "statut_name" column after "statut_label" (i.e. before the
# existing "NCTId" column), shifting all subsequent columns one to the
# right, then populate the new column with text derived from the
# statut_label value of each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; everything currently in C..L moves to D..M.
$ws.Columns.Item(3).Insert()

# Header for the new column. (The inserted column already inherits the
# correct header formatting/style from its neighbours.)
$ws.Cells.Item(1, 3).Value = "statut_name"

# Mapping from statut_label (column B) to the new statut_name text.
$map = @{
    "rouge"  = "résultat et / ou publication posté"
    "noir"   = "pas de résultat ni de publication"
    "orange" = "résultat et / ou publication posté dans les 36 mois"
    "vert"   = "résultat et / ou publication posté dans les 12 mois"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2
    if ($map.ContainsKey($label)) {
        $ws.Cells.Item($r, 3).Value = $map[$label]
    }
}
